$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.729832172393799
$ws.Range("B1").Value = 1.763488411903381
$ws.Range("C1").Value = 1.929625988006592
$ws.Range("D1").Value = 2.769926786422729
$ws.Range("E1").Value = 4.927236557006836
